# Update biofuels, calibration, heat generation
# This edits the "Output" worksheet of the Hydrogen results workbook:
#  - Row 2 (tech=HH2_BIO_SR_C_NEW) becomes tech=HH2_WE_ALK_DS_NEW,
#    output_comm HH2_CT -> HH2_WE_DT, with recalibrated 2025/2030/2035 values.
#  - Row 3 (tech=HH2_NGA_CL_CCS_NEW, output_comm HH2_CU) keeps its labels but
#    gets recalibrated 2040/2045/2050/... values.
#  - Rows 4 and 5 (HH2_COA_CL_CCS_NEW / HH2_COA_CM_CCS_NEW techs, both with
#    output_comm HH2_CT) are removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Output")

# --- Row 2: HH2_BIO_SR_C_NEW -> HH2_WE_ALK_DS_NEW, output HH2_CT -> HH2_WE_DT ---
$ws.Cells.Item(2, 2).Value = "HH2_WE_ALK_DS_NEW"   # B2 tech
$ws.Cells.Item(2, 3).Value = "HH2_WE_DT"           # C2 output_comm
$ws.Cells.Item(2, 5).Value = 0                     # E2 (2025)
$ws.Cells.Item(2, 6).Value = 0.5399999999999999    # F2 (2030)
$ws.Cells.Item(2, 7).Value = 0.02178016095559211   # G2 (2035)

# --- Row 3: HH2_NGA_CL_CCS_NEW / HH2_CU, recalibrated later years ---
$ws.Cells.Item(3, 9).Value = 1.367853537709895     # I3 (2040)
$ws.Cells.Item(3, 10).Value = 8.480691933801575    # J3 (2045)
$ws.Cells.Item(3, 11).Value = 328.9520159999997    # K3 (2050)
$ws.Cells.Item(3, 12).Value = 330.6833664732822    # L3 (2055/last col)

# --- Remove rows 4 and 5 (HH2_COA_CL_CCS_NEW, HH2_COA_CM_CCS_NEW) ---
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(4).Delete()
